$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell $ws "D2" "245.55"
Set-TextCell $ws "E2" "0.79%"
Set-TextCell $ws "D3" "29.33"
Set-TextCell $ws "E3" "-1.77%"
Set-TextCell $ws "D4" "5.143"
Set-TextCell $ws "E4" "-0.23%"
Set-TextCell $ws "D5" "0.05784"
Set-TextCell $ws "E5" "2.14%"
Set-TextCell $ws "D6" "6.627"
Set-TextCell $ws "E6" "1.58%"
Set-TextCell $ws "D7" "3.177"
Set-TextCell $ws "E7" "5.13%"
Set-TextCell $ws "D8" "0.8592"
Set-TextCell $ws "E8" "2.25%"
Set-TextCell $ws "D9" "0.8625"
Set-TextCell $ws "E9" "-0.30%"
Set-TextCell $ws "B10" "One"
Set-TextCell $ws "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws "D10" "0.01026"
Set-TextCell $ws "E10" "1.85%"
Set-TextCell $ws "B11" "WazirX"
Set-TextCell $ws "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws "D11" "0.1366"
Set-TextCell $ws "B12" "MandalaExchangeToken"
Set-TextCell $ws "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell $ws "D12" "0.07072"
Set-TextCell $ws "E12" "2.40%"
Set-TextCell $ws "B13" "BitrueCoin"
Set-TextCell $ws "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell $ws "D13" "0.03196"
Set-TextCell $ws "E13" "9.91%"
Set-TextCell $ws "B14" "BitMartToken"
Set-TextCell $ws "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell $ws "D14" "0.09354"
Set-TextCell $ws "E14" "-0.30%"
Set-TextCell $ws "B15" "BitForexToken"
Set-TextCell $ws "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws "D15" "0.001525"
Set-TextCell $ws "E15" "0.22%"
Set-TextCell $ws "D16" "0.005997"
Set-TextCell $ws "E16" "-3.27%"
Set-TextCell $ws "D17" "3.484"
Set-TextCell $ws "E17" "-0.65%"
Set-TextCell $ws "D18" "2.166"
Set-TextCell $ws "E18" "-2.91%"
Set-TextCell $ws "D20" "0.03300"
Set-TextCell $ws "E20" "1.22%"
Set-TextCell $ws "D21" "0.1284"
Set-TextCell $ws "E21" "-1.48%"
Set-TextCell $ws "D22" "3.321"
Set-TextCell $ws "E22" "-8.06%"
Set-TextCell $ws "D23" "0.04135"
Set-TextCell $ws "E23" "-0.72%"
Set-TextCell $ws "E24" "1.87%"
Set-TextCell $ws "D25" "0.001225"
Set-TextCell $ws "E25" "1.25%"
Set-TextCell $ws "D26" "0.004136"
Set-TextCell $ws "E26" "-6.89%"
Set-TextCell $ws "E27" "2.52%"
Set-TextCell $ws "E28" "3.37%"
Set-TextCell $ws "D40" "0.03733"
Set-TextCell $ws "D41" "0.005747"
Set-TextCell $ws "E41" "7.86%"
Set-TextCell $ws "D42" "0.1068"
Set-TextCell $ws "E42" "1.01%"
Set-TextCell $ws "D43" "0.001999"
Set-TextCell $ws "E43" "-13.47%"
Set-TextCell $ws "D44" "0.009179"
Set-TextCell $ws "E44" "-6.18%"
Set-TextCell $ws "D45" "0.00005266"
Set-TextCell $ws "E45" "3.41%"
Set-TextCell $ws "D46" "0.00000000750"
Set-TextCell $ws "E46" "-0.02%"
Set-TextCell $ws "D47" "0.05797"
Set-TextCell $ws "E47" "-42.01%"
Set-TextCell $ws "E48" "118.58%"
Set-TextCell $ws "D49" "0.00002099"
Set-TextCell $ws "E49" "-0.02%"
Set-TextCell $ws "D50" "0.0001999"
Set-TextCell $ws "E50" "-0.02%"

Write-Host "Applied cryptos.xlsx update: $($wb.Name)"
